$d = $word.ActiveDocument

# --- Add the new custom paragraph styles (wdStyleTypeParagraph = 1) ---

# tei_signed  (styleId "teisigned", based on Normal)
$teiSigned = $d.Styles.Add("teisigned", 1)
$teiSigned.NameLocal = "tei_signed"
$teiSigned.BaseStyle = $d.Styles("Normal")
$teiSigned.QuickStyle = $true
$teiSigned.ParagraphFormat.SpaceBefore = 18
$teiSigned.ParagraphFormat.LeftIndent = 21.55
$teiSigned.ParagraphFormat.FirstLineIndent = -21.55

# tei_speech  (styleId "teispeech", based on Normal)
$teiSpeech = $d.Styles.Add("teispeech", 1)
$teiSpeech.NameLocal = "tei_speech"
$teiSpeech.BaseStyle = $d.Styles("Normal")
$teiSpeech.QuickStyle = $true
$teiSpeech.ParagraphFormat.LeftIndent = 21.6
$teiSpeech.ParagraphFormat.FirstLineIndent = -21.6

# GeneratedTitle  (based on Title)
$genTitle = $d.Styles.Add("GeneratedTitle", 1)
$genTitle.NameLocal = "GeneratedTitle"
$genTitle.BaseStyle = $d.Styles("Title")
$genTitle.QuickStyle = $true

# GeneratedSubTitle  (based on Subtitle)
$genSubTitle = $d.Styles.Add("GeneratedSubTitle", 1)
$genSubTitle.NameLocal = "GeneratedSubTitle"
$genSubTitle.BaseStyle = $d.Styles("Subtitle")
$genSubTitle.QuickStyle = $true

# --- Re-style the lone (empty) paragraph to use the new GeneratedSubTitle
#     style instead of its explicit hanging-indent direct formatting ---
$p1 = $d.Paragraphs(1)
$p1.Style = "GeneratedSubTitle"
